$wb = $excel.ActiveWorkbook

# Each sheet gets one additional data row (row 91) appended after the
# existing last row (row 90), mirroring the date/time stamp of the new
# reading and its decoded values.

$sheetsData = @(
    @{
        Index = 1
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x2C"
        E = "0x14"
        F = 380
        G = "7.598631275147109e+23"
        H = 300
        I = 14
    },
    @{
        Index = 2
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x34"
        E = "0xe"
        F = 380
        G = "5.68432987514711e+23"
        H = 308
        I = 14
    },
    @{
        Index = 3
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x76"
        E = "0x7"
        F = 130
        G = "5.68631262647114e+23"
        H = 118
        I = 7
    },
    @{
        Index = 4
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x75"
        E = "0x3"
        F = 130
        G = "9.85046333984776e+23"
        H = 117
        I = 3
    }
)

foreach ($entry in $sheetsData) {
    $ws = $wb.Worksheets.Item($entry.Index)
    $row = 91

    $ws.Cells.Item($row, 1).Value = 45877.43872685185
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $entry.B
    $ws.Cells.Item($row, 3).Value = $entry.C
    $ws.Cells.Item($row, 4).Value = $entry.D
    $ws.Cells.Item($row, 5).Value = $entry.E
    $ws.Cells.Item($row, 6).Value = $entry.F
    $ws.Cells.Item($row, 7).Value = [double]$entry.G
    $ws.Cells.Item($row, 8).Value = $entry.H
    $ws.Cells.Item($row, 9).Value = $entry.I
}
